$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 617-618 (shifts existing rows 617-697 down to 619-699),
# matching the new weekly price update for Betarraga added to the dataset.
$ws.Rows("617:618").Insert()

# Row 617 - new "Primera" quality record for 2023-07-17
$ws.Cells.Item(617, 1).Value = 7
$ws.Cells.Item(617, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(617, 3).Value = "Ñuble"
$ws.Cells.Item(617, 4).Value = 45124
$ws.Cells.Item(617, 5).Value = 16
$ws.Cells.Item(617, 6).Value = 100114014
$ws.Cells.Item(617, 7).Value = "Betarraga"
$ws.Cells.Item(617, 8).Value = "Sin especificar"
$ws.Cells.Item(617, 9).Value = "Primera"
$ws.Cells.Item(617, 10).Value = 300
$ws.Cells.Item(617, 11).Value = 900
$ws.Cells.Item(617, 12).Value = 900
$ws.Cells.Item(617, 13).Value = 900
$ws.Cells.Item(617, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(617, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(617, 16).Value = 180
$ws.Cells.Item(617, 17).Value = 5
$ws.Cells.Item(617, 18).Value = "Hortaliza"

# Row 618 - new "Segunda" quality record for 2023-07-17
$ws.Cells.Item(618, 1).Value = 7
$ws.Cells.Item(618, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(618, 3).Value = "Ñuble"
$ws.Cells.Item(618, 4).Value = 45124
$ws.Cells.Item(618, 5).Value = 16
$ws.Cells.Item(618, 6).Value = 100114014
$ws.Cells.Item(618, 7).Value = "Betarraga"
$ws.Cells.Item(618, 8).Value = "Sin especificar"
$ws.Cells.Item(618, 9).Value = "Segunda"
$ws.Cells.Item(618, 10).Value = 250
$ws.Cells.Item(618, 11).Value = 700
$ws.Cells.Item(618, 12).Value = 700
$ws.Cells.Item(618, 13).Value = 700
$ws.Cells.Item(618, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(618, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(618, 16).Value = 140
$ws.Cells.Item(618, 17).Value = 5
$ws.Cells.Item(618, 18).Value = "Hortaliza"
